# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to the leve-profit tracking sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 579.3333
$ws.Range("I28").Value = 206.53847
$ws.Range("J28").Value = 3002.5
$ws.Range("K28").Value = 206.53847
$ws.Range("L28").Value = 3002.5
$ws.Range("M28").Value = 278.46153
$ws.Range("N28").Value = -3972.5
$ws.Range("H41").Value = 195.22223
$ws.Range("I41").Value = 151.4
$ws.Range("K41").Value = 151.4
$ws.Range("M41").Value = 288.6
$ws.Range("H62").Value = 6784.5625
$ws.Range("I62").Value = 6784.5625
$ws.Range("K62").Value = 6784.5625
$ws.Range("M62").Value = -6160.5625
$ws.Range("H64").Value = 4424.4546
$ws.Range("I64").Value = 3600
$ws.Range("J64").Value = 4895.5713
$ws.Range("K64").Value = 3600
$ws.Range("L64").Value = 4895.5713
$ws.Range("M64").Value = -3352
$ws.Range("N64").Value = -5391.5713
$ws.Range("H65").Value = 6784.5625
$ws.Range("I65").Value = 6784.5625
$ws.Range("K65").Value = 33922.8125
$ws.Range("M65").Value = -30802.8125
$ws.Range("H67").Value = 4424.4546
$ws.Range("I67").Value = 3600
$ws.Range("J67").Value = 4895.5713
$ws.Range("K67").Value = 3600
$ws.Range("L67").Value = 4895.5713
$ws.Range("M67").Value = -2742
$ws.Range("N67").Value = -6611.5713
$ws.Range("H99").Value = 563
$ws.Range("I99").Value = 563
$ws.Range("K99").Value = 1689
$ws.Range("M99").Value = -191
$ws.Range("H100").Value = 49497.137
$ws.Range("I100").Value = 59218.555
$ws.Range("K100").Value = 59218.555
$ws.Range("M100").Value = -58677.555
$ws.Range("H112").Value = 2446.9429
$ws.Range("J112").Value = 2490.1765
$ws.Range("L112").Value = 7470.529500000001
$ws.Range("N112").Value = -9686.529500000001
$ws.Range("H132").Value = 1364.3243
$ws.Range("I132").Value = 1368.4062
$ws.Range("K132").Value = 4105.2186
$ws.Range("M132").Value = -1575.2186
$ws.Range("H137").Value = 9397.811
$ws.Range("I137").Value = 4946.4443
$ws.Range("J137").Value = 13614.895
$ws.Range("K137").Value = 14839.3329
$ws.Range("L137").Value = 40844.685
$ws.Range("M137").Value = -12289.3329
$ws.Range("N137").Value = -45944.685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5796.2173
$ws.Range("I61").Value = 3818.3928
$ws.Range("J61").Value = 8872.833000000001
$ws.Range("K61").Value = 3818.3928
$ws.Range("L61").Value = 8872.833000000001
$ws.Range("M61").Value = -3606.3928
$ws.Range("N61").Value = -9296.833000000001
$ws.Range("H74").Value = 19418.818
$ws.Range("J74").Value = 7083.3335
$ws.Range("L74").Value = 7083.3335
$ws.Range("N74").Value = -8831.333500000001
$ws.Range("H77").Value = 19418.818
$ws.Range("J77").Value = 7083.3335
$ws.Range("L77").Value = 35416.6675
$ws.Range("N77").Value = -44152.6675
$ws.Range("H102").Value = 1866.7646
$ws.Range("I102").Value = 1937.25
$ws.Range("J102").Value = 1697.6
$ws.Range("K102").Value = 1937.25
$ws.Range("L102").Value = 1697.6
$ws.Range("M102").Value = -315.25
$ws.Range("N102").Value = -4941.6
$ws.Range("H130").Value = 63142
$ws.Range("J130").Value = 63142
$ws.Range("L130").Value = 63142
$ws.Range("N130").Value = -73182
$ws.Range("H136").Value = 5796.2173
$ws.Range("I136").Value = 3818.3928
$ws.Range("J136").Value = 8872.833000000001
$ws.Range("K136").Value = 11455.1784
$ws.Range("L136").Value = 26618.499
$ws.Range("M136").Value = -8905.178400000001
$ws.Range("N136").Value = -31718.499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9954.058000000001
$ws.Range("I20").Value = 12766.462
$ws.Range("J20").Value = 1829.3334
$ws.Range("K20").Value = 12766.462
$ws.Range("L20").Value = 1829.3334
$ws.Range("M20").Value = -12519.462
$ws.Range("N20").Value = -2323.3334
$ws.Range("H103").Value = 34067.5
$ws.Range("J103").Value = 34067.5
$ws.Range("L103").Value = 34067.5
$ws.Range("N103").Value = -36411.5
$ws.Range("H105").Value = 5874.3335
$ws.Range("I105").Value = 5686.875
$ws.Range("J105").Value = 6249.25
$ws.Range("K105").Value = 5686.875
$ws.Range("L105").Value = 6249.25
$ws.Range("M105").Value = -3939.875
$ws.Range("N105").Value = -9743.25
$ws.Range("H109").Value = 40000
$ws.Range("J109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("N109").Value = -42774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 37020
$ws.Range("J92").Value = 37020
$ws.Range("L92").Value = 37020
$ws.Range("N92").Value = -42012
$ws.Range("H132").Value = 29279.582
$ws.Range("J132").Value = 34798.125
$ws.Range("L132").Value = 104394.375
$ws.Range("N132").Value = -109454.375
$ws.Range("H134").Value = 4394.4707
$ws.Range("I134").Value = 2635.1724
$ws.Range("J134").Value = 14598.4
$ws.Range("K134").Value = 7905.5172
$ws.Range("L134").Value = 43795.2
$ws.Range("M134").Value = -5370.5172
$ws.Range("N134").Value = -48865.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H5").Value = 1712.0344
$ws.Range("I5").Value = 820.7778
$ws.Range("J5").Value = 2113.1
$ws.Range("K5").Value = 2462.3334
$ws.Range("L5").Value = 6339.299999999999
$ws.Range("M5").Value = -2350.3334
$ws.Range("N5").Value = -6563.299999999999
$ws.Range("H129").Value = 2276.4443
$ws.Range("J129").Value = 2727.1428
$ws.Range("L129").Value = 8181.428400000001
$ws.Range("N129").Value = -18181.4284
$ws.Range("H133").Value = 6391
$ws.Range("J133").Value = 6391
$ws.Range("L133").Value = 19173
$ws.Range("N133").Value = -29293
$ws.Range("H135").Value = 1712.0344
$ws.Range("I135").Value = 820.7778
$ws.Range("J135").Value = 2113.1
$ws.Range("K135").Value = 7387.000199999999
$ws.Range("L135").Value = 19017.9
$ws.Range("M135").Value = -4852.000199999999
$ws.Range("N135").Value = -24087.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7774.625
$ws.Range("I70").Value = 8268.333000000001
$ws.Range("K70").Value = 8268.333000000001
$ws.Range("M70").Value = -7998.333000000001
$ws.Range("H73").Value = 7774.625
$ws.Range("I73").Value = 8268.333000000001
$ws.Range("K73").Value = 8268.333000000001
$ws.Range("M73").Value = -7332.333000000001
$ws.Range("H97").Value = 483.70587
$ws.Range("I97").Value = 505.33334
$ws.Range("J97").Value = 459.375
$ws.Range("K97").Value = 505.33334
$ws.Range("L97").Value = 459.375
$ws.Range("M97").Value = -9.333340000000021
$ws.Range("N97").Value = -1451.375
$ws.Range("H132").Value = 16746.87
$ws.Range("I132").Value = 12220.857
$ws.Range("K132").Value = 36662.571
$ws.Range("M132").Value = -34132.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8681.083000000001
$ws.Range("I132").Value = 8020.4443
$ws.Range("J132").Value = 10663
$ws.Range("K132").Value = 24061.3329
$ws.Range("L132").Value = 31989
$ws.Range("M132").Value = -21531.3329
$ws.Range("N132").Value = -37049
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 4000.175
$ws.Range("I136").Value = 3534.3076
$ws.Range("K136").Value = 10602.9228
$ws.Range("M136").Value = -8052.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5292646.5
$ws.Range("I107").Value = 1112.5
$ws.Range("J107").Value = 15875714
$ws.Range("K107").Value = 3337.5
$ws.Range("L107").Value = 47627142
$ws.Range("M107").Value = -1417.5
$ws.Range("N107").Value = -47630982
$ws.Range("H113").Value = 690.6316
$ws.Range("I113").Value = 366.85715
$ws.Range("K113").Value = 1100.57145
$ws.Range("M113").Value = 1069.42855
$ws.Range("H132").Value = 155173.53
$ws.Range("I132").Value = 298079.8
$ws.Range("J132").Value = 19991.918
$ws.Range("K132").Value = 894239.3999999999
$ws.Range("L132").Value = 59975.754
$ws.Range("M132").Value = -891709.3999999999
$ws.Range("N132").Value = -65035.754
$ws.Range("H136").Value = 3775776.8
$ws.Range("I136").Value = 6668201.5
$ws.Range("K136").Value = 20004604.5
$ws.Range("M136").Value = -20002054.5
$ws.Range("H140").Value = 64332.668
$ws.Range("J140").Value = 64332.668
$ws.Range("L140").Value = 64332.668
$ws.Range("N140").Value = -74692.66800000001
